# Update the cryptocurrency price/volume table with the latest scraped
# values. A leading apostrophe is used on numeric-looking Price values so
# Excel stores them as literal text (matching the original formatting,
# e.g. "64.37") instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "41.730.08"
$ws.Range("E2").Value = "  +5.77%  "
$ws.Range("D3").Value = "2.257.31"
$ws.Range("E3").Value = "  +4.39%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'233.16"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("D7").Value = "'64.37"
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.415"
$ws.Range("E9").Value = "  +5.07%  "
$ws.Range("D10").Value = "'60.27"
$ws.Range("E10").Value = "  +3.76%  "
$ws.Range("D11").Value = "'0.0909"
$ws.Range("E11").Value = "  +6.41%  "
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").Value = "2.590.41"
$ws.Range("E13").Value = "  +4.27%  "
$ws.Range("D14").Value = "'16.25"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "'22.71"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").Value = "2.252.46"
$ws.Range("E18").Value = "  +4.01%  "
$ws.Range("D19").Value = "41.562.47"
$ws.Range("E19").Value = "  +5.45%  "
$ws.Range("E20").Value = "  +9.56%  "
$ws.Range("D21").Value = "'74.05"
$ws.Range("E21").Value = "  +2.97%  "
$ws.Range("D22").Value = "'6.21"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D23").Value = "'253.51"
$ws.Range("E23").Value = "  +10.58%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'2.40"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "'2.35"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.150"
$ws.Range("E27").Value = "  +5.84%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'9.87"
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").Value = "'172.08"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'20.54"
$ws.Range("E30").Value = "  +3.32%  "
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").Value = "'2.82"
$ws.Range("E32").Value = "  +8.15%  "
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").Value = "'5.13"
$ws.Range("E34").Value = "  +8.79%  "
$ws.Range("D35").Value = "'4.78"
$ws.Range("E35").Value = "  +3.47%  "
$ws.Range("D36").Value = "'0.0643"
$ws.Range("E36").Value = "  +4.09%  "
$ws.Range("D37").Value = "'6.93"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").Value = "'3.85"
$ws.Range("E38").Value = "  +7.83%  "
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("D40").Value = "'0.000259"
$ws.Range("E40").Value = "  +63.76%  "
$ws.Range("D41").Value = "'5.14"
$ws.Range("E41").Value = "  +20.37%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").Value = "'0.0244"
$ws.Range("E43").Value = "  +6.22%  "
$ws.Range("D44").Value = "'8.82"
$ws.Range("E44").Value = "  +13.97%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'102.90"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.1000"
$ws.Range("E46").Value = "  +7.52%  "
$ws.Range("D47").Value = "'17.71"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "'1.23"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("D49").Value = "1.512.92"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "'1.15"
$ws.Range("E50").Value = "  +3.65%  "
$ws.Range("E51").Value = "  -0.91%  "
